# Daily attendance processing - 2025-12-31 14:59:05
# Swap the "Recorded By" column (G) entries that list "System" first so
# that the human reviewer's email appears before "System", e.g.
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
